$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.899.94"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.41%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.744.47"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.31%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "231.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.97%  "
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5176"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.98%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2805"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.38%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "39.53"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06116"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.14%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.752.90"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.20%  "
$ws.Range("E12").Value = "  +1.61%  "
$ws.Range("E13").Value = "  -0.79%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6411"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.26%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.519"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.99%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "77.21"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.65%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.001"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.03%  "
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "25.881.81"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.30%  "
$ws.Range("E20").Value = "  -1.32%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000006584"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.48%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.973.16"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.140"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.633"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.61%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.141"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.25%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "140.22"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.60%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.515"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.29%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.11"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("E29").Value = "  +3.86%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "102.25"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.21%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08228"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.49%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.667"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.75%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.427"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.83%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04488"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.52%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.615"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.32%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9854"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6142"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.54%  "
$ws.Range("E38").Value = "  +1.11%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01592"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.77%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.923"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.25%  "
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "100.74"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.49%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.3842"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.52%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.028"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.88%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.7221"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.86%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.05430"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.31%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.282"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.64%  "
$ws.Range("E48").Value = "  +2.20%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "53.07"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.664"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.78%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "29.86"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.71%  "
